# ITRC-ETR-2025-0001.docx content refresh: placeholder/demo text -> final
# report values ("true report generation with help as v3").
#
# Implementation notes (both verified experimentally against this runtime):
#
#  1. Range.Find.Execute always searches from the start of the active story
#     and ignores the calling Range's own Start/End, so a Find/Replace
#     "scoped" to one cell/paragraph can silently hit the *first* match
#     anywhere in the document (e.g. it would turn the untouched
#     "2025-06-02 to 2025-06-02" Evaluation Period cell into
#     "2025-06-03 to 2025-06-03" when we only meant to touch the Report
#     Date cell). So plain Find/Replace is avoided below.
#
#  2. $d.Paragraphs is a cached collection that gets thrown out of sync
#     (its Item(n) starts resolving against the wrong range) as soon as a
#     Tables/Cell access happens in between. Re-deriving the paragraph
#     collection from $d.Content each time (`$d.Content.Paragraphs`) keeps
#     indices stable and absolute no matter what table access happened
#     before it.
#
# Every paragraph/cell touched here holds exactly one run with no mixed
# formatting, so assigning Range.Text directly retargets each value
# precisely, leaves w:rPr (bold/size) untouched, and never disturbs
# look-alike text elsewhere (e.g. "dfgdfgdfg" containing "dfgdfgdf",
# "FAM_MAL.1" containing "FAM_MAL").

$d = $word.ActiveDocument

# --- Title / subtitle (paragraphs 3 and 5) ---------------------------------
$d.Content.Paragraphs.Item(3).Range.Text = "Evaluation Technical Report for this is test"
$d.Content.Paragraphs.Item(5).Range.Text = "Product: this is test"

# --- Summary table (Table 1) ------------------------------------------------
$d.Tables.Item(1).Cell(1, 2).Range.Text = "ITRC-2025-1748854557"
$d.Tables.Item(1).Cell(4, 2).Range.Text = "dfdf"
$d.Tables.Item(1).Cell(6, 2).Range.Text = "2025-06-03"

# --- Executive Summary paragraph -------------------------------------------
$d.Content.Paragraphs.Item(34).Range.Text = "This report presents the results of the Common Criteria evaluation of this is test version N/A conducted by the Iran Telecommunications Research Center (ITRC)."

# --- Product Identification table (Table 2) ---------------------------------
$d.Tables.Item(2).Cell(1, 2).Range.Text = "this is test"
$d.Tables.Item(2).Cell(3, 2).Range.Text = "dfdf"
$d.Tables.Item(2).Cell(5, 2).Range.Text = "ITRC-2025-1748854557"

# --- Product Description -----------------------------------------------------
$d.Content.Paragraphs.Item(54).Range.Text = "sdfsdf"

# --- Section heading "1. Malware Detection" -> "1. Scanning Engine" --------
$d.Content.Paragraphs.Item(66).Range.Text = "1. Scanning Engine"

# --- Detailed results table (Table 3) ---------------------------------------
$d.Tables.Item(3).Cell(1, 2).Range.Text = "FAM_SCN"
$d.Tables.Item(3).Cell(2, 2).Range.Text = "On-demand scanning (FAM_SCN.1)"
$d.Tables.Item(3).Cell(4, 2).Range.Text = "N/A"
$d.Tables.Item(3).Cell(5, 2).Range.Text = "موتور اسکن"

# --- Implementation / Justification / Test Approach / Evaluator Assessment -
$d.Content.Paragraphs.Item(83).Range.Text = "dfdf"
$d.Content.Paragraphs.Item(85).Range.Text = "sfg"
$d.Content.Paragraphs.Item(87).Range.Text = "hrhr"
$d.Content.Paragraphs.Item(89).Range.Text = "rerer"
